$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 41249014.21622031
$ws.Range("D2").Value = 189.6080260415259
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 41249220.96042006

# Row 3 updates
$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 3.223369029078222
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 6.82939032824165
